$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48-107 down to 49-108.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new data record.
$ws.Cells.Item(48, 1).Value() = 4
$ws.Cells.Item(48, 2).Value() = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(48, 3).Value() = "Los Lagos"
$ws.Cells.Item(48, 4).Value() = 44467
$ws.Cells.Item(48, 5).Value() = 10
$ws.Cells.Item(48, 6).Value() = 100112009
$ws.Cells.Item(48, 7).Value() = "Acelga"
$ws.Cells.Item(48, 8).Value() = "Sin especificar"
$ws.Cells.Item(48, 9).Value() = "Primera"
$ws.Cells.Item(48, 10).Value() = 200
$ws.Cells.Item(48, 11).Value() = 4000
$ws.Cells.Item(48, 12).Value() = 4000
$ws.Cells.Item(48, 13).Value() = 4000
$ws.Cells.Item(48, 14).Value() = "$/docena de atados (4 kilos)"
$ws.Cells.Item(48, 15).Value() = "Región del Maule"
$ws.Cells.Item(48, 16).Value() = 1000
$ws.Cells.Item(48, 17).Value() = 4
$ws.Cells.Item(48, 18).Value() = "Hortaliza"
